$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row "A" plate readings (B34:M34) were cleared out -- keeping the cell
# styling but dropping the measured values (supports lowercase row
# indexing going forward, per commit message).
$ws.Range("B34:M34").ClearContents()

# Reflect the user's on-screen state at the time of the edit: scrolled
# down to row 23 with the now-empty B34:M34 range selected.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 23
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B34:M34").Select()
